$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.295.37"
$ws.Range("E2").Value = "  -2.43%  "
$ws.Range("D3").Value = "'1.870.20"
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'319.19"
$ws.Range("E5").Value = "  -1.66%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D7").Value = "'0.4421"
$ws.Range("E7").Value = "  -3.83%  "
$ws.Range("D8").Value = "'0.3700"
$ws.Range("E8").Value = "  -3.28%  "
$ws.Range("D9").Value = "'0.07518"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("D10").Value = "'0.9394"
$ws.Range("E10").Value = "  -4.17%  "
$ws.Range("D11").Value = "'21.44"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "'1.907.28"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'6.716"
$ws.Range("E13").Value = "  -3.13%  "
$ws.Range("D14").Value = "'5.469"
$ws.Range("E14").Value = "  -3.47%  "
$ws.Range("D15").Value = "'0.06882"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "'1.003"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "'82.18"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "'0.000009059"
$ws.Range("E18").Value = "  -4.37%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("D20").Value = "'15.95"
$ws.Range("E20").Value = "  -4.35%  "
$ws.Range("D21").Value = "'28.287.18"
$ws.Range("E21").Value = "  -2.41%  "
$ws.Range("D22").Value = "'5.129"
$ws.Range("E22").Value = "  -3.51%  "
$ws.Range("D23").Value = "'10.87"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'2.133.52"
$ws.Range("E24").Value = "  -0.82%  "
$ws.Range("D25").Value = "'2.025"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("D26").Value = "'154.95"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").Value = "'18.43"
$ws.Range("E27").Value = "  -3.20%  "
$ws.Range("D28").Value = "'5.346"
$ws.Range("E28").Value = "  -5.40%  "
$ws.Range("D29").Value = "'113.84"
$ws.Range("E29").Value = "  -2.98%  "
$ws.Range("E30").Value = "  -6.79%  "
$ws.Range("D31").Value = "'0.09065"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("D32").Value = "'0.8007"
$ws.Range("E32").Value = "  -7.05%  "
$ws.Range("D33").Value = "'4.867"
$ws.Range("E33").Value = "  -4.16%  "
$ws.Range("D34").Value = "'1.179"
$ws.Range("E34").Value = "  -5.22%  "
$ws.Range("D35").Value = "'2.919"
$ws.Range("E35").Value = "  -1.17%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'1.128"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "'0.05451"
$ws.Range("E38").Value = "  -4.85%  "
$ws.Range("D39").Value = "'0.01971"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").Value = "'3.010"
$ws.Range("E40").Value = "  +6.46%  "
$ws.Range("D41").Value = "'7.128"
$ws.Range("E41").Value = "  -3.70%  "
$ws.Range("D42").Value = "'0.5265"
$ws.Range("E42").Value = "  -4.33%  "
$ws.Range("D43").Value = "'0.1688"
$ws.Range("E43").Value = "  -3.92%  "
$ws.Range("E44").Value = "  -5.88%  "
$ws.Range("D45").Value = "'0.06761"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").Value = "'0.4890"
$ws.Range("E46").Value = "  -5.63%  "
$ws.Range("D47").Value = "'10.63"
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("D48").Value = "'107.83"
$ws.Range("E48").Value = "  -2.87%  "
$ws.Range("D49").Value = "'1.962"
$ws.Range("E49").Value = "  -4.49%  "
$ws.Range("D50").Value = "'0.000002441"
$ws.Range("E50").Value = "  -3.32%  "
$ws.Range("D51").Value = "'1.680"
$ws.Range("E51").Value = "  -5.36%  "
